$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Rows 11 and 12: coin name/link swap plus price/volume updates
Set-TextValue $ws.Range("B11") "MandalaExchangeToken"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.09602"
Set-TextValue $ws.Range("E11") "3.05%"
Set-TextValue $ws.Range("B12") "BitrueCoin"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.04732"
Set-TextValue $ws.Range("E12") "2.91%"

# Price (D) and Volume(1h) (E) updates
Set-TextValue $ws.Range("D2") "309.94"
Set-TextValue $ws.Range("E2") "-2.64%"
Set-TextValue $ws.Range("D3") "52.28"
Set-TextValue $ws.Range("E3") "7.75%"
Set-TextValue $ws.Range("D4") "5.121"
Set-TextValue $ws.Range("E4") "-2.75%"
Set-TextValue $ws.Range("D5") "0.07796"
Set-TextValue $ws.Range("E5") "-2.41%"
Set-TextValue $ws.Range("D6") "4.495"
Set-TextValue $ws.Range("E6") "-2.11%"
Set-TextValue $ws.Range("D7") "1.356"
Set-TextValue $ws.Range("E7") "-5.35%"
Set-TextValue $ws.Range("D8") "1.582"
Set-TextValue $ws.Range("E8") "-3.99%"
Set-TextValue $ws.Range("D9") "0.1227"
Set-TextValue $ws.Range("E9") "-3.95%"
Set-TextValue $ws.Range("E10") "3.49%"
Set-TextValue $ws.Range("E13") "0.05%"
Set-TextValue $ws.Range("D14") "0.001263"
Set-TextValue $ws.Range("E14") "-4.22%"
Set-TextValue $ws.Range("D15") "0.005791"
Set-TextValue $ws.Range("E15") "-0.56%"
Set-TextValue $ws.Range("E16") "2,013.74%"
Set-TextValue $ws.Range("D17") "3.334"
Set-TextValue $ws.Range("E17") "0.02%"
Set-TextValue $ws.Range("D18") "2.414"
Set-TextValue $ws.Range("E18") "-0.86%"
Set-TextValue $ws.Range("E19") "1.60%"
Set-TextValue $ws.Range("D20") "8.024"
Set-TextValue $ws.Range("E20") "-1.79%"
Set-TextValue $ws.Range("D21") "0.1364"
Set-TextValue $ws.Range("E21") "-2.66%"
Set-TextValue $ws.Range("E22") "-0.23%"
Set-TextValue $ws.Range("D23") "0.04174"
Set-TextValue $ws.Range("E23") "-0.03%"
Set-TextValue $ws.Range("D24") "0.001262"
Set-TextValue $ws.Range("E24") "-3.55%"
Set-TextValue $ws.Range("D25") "0.003945"
Set-TextValue $ws.Range("E25") "-6.92%"
Set-TextValue $ws.Range("D26") "0.0001352"
Set-TextValue $ws.Range("E26") "0.07%"
Set-TextValue $ws.Range("E38") "-3.94%"
Set-TextValue $ws.Range("E39") "3.89%"
Set-TextValue $ws.Range("D40") "0.01101"
Set-TextValue $ws.Range("E40") "74.38%"
Set-TextValue $ws.Range("D41") "0.008082"
Set-TextValue $ws.Range("E41") "0.98%"
Set-TextValue $ws.Range("D42") "0.1423"
Set-TextValue $ws.Range("E42") "-1.16%"
Set-TextValue $ws.Range("D43") "0.008244"
Set-TextValue $ws.Range("E43") "7.27%"
Set-TextValue $ws.Range("D44") "0.008443"
Set-TextValue $ws.Range("E44") "7.08%"
Set-TextValue $ws.Range("D45") "0.3118"
Set-TextValue $ws.Range("E45") "-10.40%"
Set-TextValue $ws.Range("D46") "0.00007352"
Set-TextValue $ws.Range("E46") "6.49%"
Set-TextValue $ws.Range("E47") "0.05%"
Set-TextValue $ws.Range("D48") "0.05756"
Set-TextValue $ws.Range("E48") "4.93%"
Set-TextValue $ws.Range("E49") "-34.48%"
Set-TextValue $ws.Range("E50") "0.05%"
Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "0.05%"
